$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new row of data: Sr No 2, Detail "This is Version 2 of File"
$ws.Range("B6").Value = 2
$ws.Range("C6").Value = "This is Version 2 of File"

# Set column C width to match the authored width (32.140625 chars).
# The COM ColumnWidth setter quantizes to the engine's pixel grid, so 31.25
# is the input that lands on the closest representable stored width.
$ws.Columns.Item(3).ColumnWidth = 31.25
